# Update the "address" (lat/lng-like) values in column D of the "Child"
# sheet (rows 2-22) to match the new set of coordinates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Child")

$values = @{
    2  = "-0.75,-7.75"
    3  = "-7.45,3.53"
    4  = "-1.07,-9.07"
    5  = "9.12,0.07"
    6  = "8.52,-5.98"
    7  = "9.35,-5.18"
    8  = "8.06,7.39"
    9  = "-7.84,3.24"
    10 = "1.1,-7.16"
    11 = "-6.44,9.6"
    12 = "-3.45,-0.28"
    13 = "1.09,-0.75"
    14 = "-2.25,-1.67"
    15 = "-8.3,6.81"
    16 = "5.68,-4.32"
    17 = "-3.23,2.78"
    18 = "9.04,-9.54"
    19 = "-5.94,3.44"
    20 = "8.28,-3.72"
    21 = "-2.97,1.58"
    22 = "6.33,5.28"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 4).Value = $values[$row]
}
